$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (this also updates the _FilterDatabase defined name
# reference automatically since it points at this sheet)
$ws.Name = "output results"

# Row 4 (F4:AK4) should share the same (non-wrapping) style as rows 3 and 5
# instead of the wrap-text style it currently has
$ws.Range("F4:AK4").WrapText = $false
